$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-01 06:35:13"

for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
